# Auto-generated Excel COM-interop script to update TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 140.209918
$ws.Range("H2").Value = 420.629754
$ws.Range("I2").Value = 0.2353423205412711
$ws.Range("J2").Value = 0.2363790708159033
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.506715
$ws.Range("N2").Value = 1.520145
$ws.Range("O2").Value = 0.003122343715987576
$ws.Range("P2").Value = 0.003132472094339857
$ws.Range("Q2").Value = 71.04646859937
$ws.Range("R2").Value = 639.41821739433
$ws.Range("S2").Value = 0.0007348196156479716
$ws.Range("T2").Value = 0.0007404508430168019
$ws.Range("G3").Value = 140.209918
$ws.Range("H3").Value = 420.629754
$ws.Range("I3").Value = 0.2353423205412711
$ws.Range("J3").Value = 0.2363790708159033
$ws.Range("M3").Value = 88.13219433333332
$ws.Range("N3").Value = 264.396583
$ws.Range("O3").Value = 0.5430646480820168
$ws.Range("P3").Value = 0.5448262620252092
$ws.Range("Q3").Value = 12357.00774063673
$ws.Range("R3").Value = 111213.0696657306
$ws.Range("S3").Value = 0.1278060944835506
$ws.Range("T3").Value = 0.1287855255736208
$ws.Range("G4").Value = 140.209918
$ws.Range("H4").Value = 420.629754
$ws.Range("I4").Value = 0.2353423205412711
$ws.Range("J4").Value = 0.2363790708159033
$ws.Range("M4").Value = 1.5741895
$ws.Range("N4").Value = 3.148379
$ws.Range("O4").Value = 0.009700049718478087
$ws.Range("P4").Value = 0.006487676741301404
$ws.Range("Q4").Value = 220.716980711461
$ws.Range("R4").Value = 1324.301884268766
$ws.Range("S4").Value = 0.002282832210112336
$ws.Range("T4").Value = 0.001533550999862773
$ws.Range("G5").Value = 140.209918
$ws.Range("H5").Value = 420.629754
$ws.Range("I5").Value = 0.2353423205412711
$ws.Range("J5").Value = 0.2363790708159033
$ws.Range("M5").Value = 72.07364666666666
$ws.Range("N5").Value = 216.22094
$ws.Range("O5").Value = 0.4441129584835175
$ws.Range("P5").Value = 0.4455535891391496
$ws.Range("Q5").Value = 10105.44008909431
$ws.Range("R5").Value = 90948.96080184875
$ws.Range("S5").Value = 0.1045185742319602
$ws.Range("T5").Value = 0.1053195433994029
$ws.Range("I6").Value = 0.7497327998952026
$ws.Range("J6").Value = 0.753035587444864
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.506715
$ws.Range("N6").Value = 1.520145
$ws.Range("O6").Value = 0.003122343715987576
$ws.Range("P6").Value = 0.003132472094339857
$ws.Range("Q6").Value = 226.333571043085
$ws.Range("R6").Value = 2037.002139387765
$ws.Range("S6").Value = 0.002340923496422557
$ws.Range("T6").Value = 0.002358862963715858
$ws.Range("I7").Value = 0.7497327998952026
$ws.Range("J7").Value = 0.753035587444864
$ws.Range("M7").Value = 88.13219433333332
$ws.Range("N7").Value = 264.396583
$ws.Range("O7").Value = 0.5430646480820168
$ws.Range("P7").Value = 0.5448262620252092
$ws.Range("Q7").Value = 39365.86496813094
$ws.Range("R7").Value = 354292.7847131785
$ws.Range("S7").Value = 0.4071533791306333
$ws.Range("T7").Value = 0.4102735642795428
$ws.Range("I8").Value = 0.7497327998952026
$ws.Range("J8").Value = 0.753035587444864
$ws.Range("M8").Value = 1.5741895
$ws.Range("N8").Value = 3.148379
$ws.Range("O8").Value = 0.009700049718478087
$ws.Range("P8").Value = 0.006487676741301404
$ws.Range("Q8").Value = 703.1406826984171
$ws.Range("R8").Value = 4218.844096190503
$ws.Range("S8").Value = 0.007272445434557248
$ws.Range("T8").Value = 0.004885451466038284
$ws.Range("I9").Value = 0.7497327998952026
$ws.Range("J9").Value = 0.753035587444864
$ws.Range("M9").Value = 72.07364666666666
$ws.Range("N9").Value = 216.22094
$ws.Range("O9").Value = 0.4441129584835175
$ws.Range("P9").Value = 0.4455535891391496
$ws.Range("Q9").Value = 32193.01940570973
$ws.Range("R9").Value = 289737.1746513876
$ws.Range("S9").Value = 0.3329660518335894
$ws.Range("T9").Value = 0.3355177087355671
$ws.Range("G10").Value = 0.7424606666666667
$ws.Range("H10").Value = 2.227382
$ws.Range("I10").Value = 0.001246220086969543
$ws.Range("J10").Value = 0.001251710043108525
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.506715
$ws.Range("N10").Value = 1.520145
$ws.Range("O10").Value = 0.003122343715987576
$ws.Range("P10").Value = 0.003132472094339857
$ws.Range("Q10").Value = 0.37621595671
$ws.Range("R10").Value = 3.38594361039
$ws.Range("S10").Value = 0.000003891127457286843
$ws.Range("T10").Value = 0.000003920946780242394
$ws.Range("G11").Value = 0.7424606666666667
$ws.Range("H11").Value = 2.227382
$ws.Range("I11").Value = 0.001246220086969543
$ws.Range("J11").Value = 0.001251710043108525
$ws.Range("M11").Value = 88.13219433333332
$ws.Range("N11").Value = 264.396583
$ws.Range("O11").Value = 0.5430646480820168
$ws.Range("P11").Value = 0.5448262620252092
$ws.Range("Q11").Value = 65.43468775952287
$ws.Range("R11").Value = 588.9121898357059
$ws.Range("S11").Value = 0.0006767780729628553
$ws.Range("T11").Value = 0.0006819645039262312
$ws.Range("G12").Value = 0.7424606666666667
$ws.Range("H12").Value = 2.227382
$ws.Range("I12").Value = 0.001246220086969543
$ws.Range("J12").Value = 0.001251710043108525
$ws.Range("M12").Value = 1.5741895
$ws.Range("N12").Value = 3.148379
$ws.Range("O12").Value = 0.009700049718478087
$ws.Range("P12").Value = 0.006487676741301404
$ws.Range("Q12").Value = 1.168773785629666
$ws.Range("R12").Value = 7.012642713778
$ws.Range("S12").Value = 0.00001208839680377065
$ws.Range("T12").Value = 0.000008120690133528557
$ws.Range("G13").Value = 0.7424606666666667
$ws.Range("H13").Value = 2.227382
$ws.Range("I13").Value = 0.001246220086969543
$ws.Range("J13").Value = 0.001251710043108525
$ws.Range("M13").Value = 72.07364666666666
$ws.Range("N13").Value = 216.22094
$ws.Range("O13").Value = 0.4441129584835175
$ws.Range("P13").Value = 0.4455535891391496
$ws.Range("Q13").Value = 53.51184775323111
$ws.Range("R13").Value = 481.60662977908
$ws.Range("S13").Value = 0.0005534624897456301
$ws.Range("T13").Value = 0.0005577039022685229
$ws.Range("G14").Value = 7.839080000000001
$ws.Range("H14").Value = 15.67816
$ws.Range("I14").Value = 0.01315789428040795
$ws.Range("J14").Value = 0.008810572380248361
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.506715
$ws.Range("N14").Value = 1.520145
$ws.Range("O14").Value = 0.003122343715987576
$ws.Range("P14").Value = 0.003132472094339857
$ws.Range("Q14").Value = 3.972179422200001
$ws.Range("R14").Value = 23.8330765332
$ws.Range("S14").Value = 0.00004108346852206063
$ws.Range("T14").Value = 0.00002759887211628948
$ws.Range("G15").Value = 7.839080000000001
$ws.Range("H15").Value = 15.67816
$ws.Range("I15").Value = 0.01315789428040795
$ws.Range("J15").Value = 0.008810572380248361
$ws.Range("M15").Value = 88.13219433333332
$ws.Range("N15").Value = 264.396583
$ws.Range("O15").Value = 0.5430646480820168
$ws.Range("P15").Value = 0.5448262620252092
$ws.Range("Q15").Value = 690.8753219545466
$ws.Range("R15").Value = 4145.25193172728
$ws.Range("S15").Value = 0.007145587226890126
$ws.Range("T15").Value = 0.004800231216233265
$ws.Range("G16").Value = 7.839080000000001
$ws.Range("H16").Value = 15.67816
$ws.Range("I16").Value = 0.01315789428040795
$ws.Range("J16").Value = 0.008810572380248361
$ws.Range("M16").Value = 1.5741895
$ws.Range("N16").Value = 3.148379
$ws.Range("O16").Value = 0.009700049718478087
$ws.Range("P16").Value = 0.006487676741301404
$ws.Range("Q16").Value = 12.34019742566
$ws.Range("R16").Value = 49.36078970264001
$ws.Range("S16").Value = 0.0001276322287104356
$ws.Range("T16").Value = 0.00005716014550888984
$ws.Range("G17").Value = 7.839080000000001
$ws.Range("H17").Value = 15.67816
$ws.Range("I17").Value = 0.01315789428040795
$ws.Range("J17").Value = 0.008810572380248361
$ws.Range("M17").Value = 72.07364666666666
$ws.Range("N17").Value = 216.22094
$ws.Range("O17").Value = 0.4441129584835175
$ws.Range("P17").Value = 0.4455535891391496
$ws.Range("Q17").Value = 564.9910821117334
$ws.Range("R17").Value = 3389.9464926704
$ws.Range("S17").Value = 0.005843591356285328
$ws.Range("T17").Value = 0.003925582146389917
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.3102563333333333
$ws.Range("H18").Value = 0.930769
$ws.Range("I18").Value = 0.0005207651961489113
$ws.Range("J18").Value = 0.0005230593158758034
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0.506715
$ws.Range("N18").Value = 1.520145
$ws.Range("O18").Value = 0.003122343715987576
$ws.Range("P18").Value = 0.003132472094339857
$ws.Range("Q18").Value = 0.157211537945
$ws.Range("R18").Value = 1.414903841505
$ws.Range("S18").Value = 0.000001626007937700591
$ws.Range("T18").Value = 0.000001638468710665451
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.3102563333333333
$ws.Range("H19").Value = 0.930769
$ws.Range("I19").Value = 0.0005207651961489113
$ws.Range("J19").Value = 0.0005230593158758034
$ws.Range("M19").Value = 88.13219433333332
$ws.Range("N19").Value = 264.396583
$ws.Range("O19").Value = 0.5430646480820168
$ws.Range("P19").Value = 0.5448262620252092
$ws.Range("Q19").Value = 27.34357146248077
$ws.Range("R19").Value = 246.0921431623269
$ws.Range("S19").Value = 0.000282809167979971
$ws.Range("T19").Value = 0.0002849764518860771
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.3102563333333333
$ws.Range("H20").Value = 0.930769
$ws.Range("I20").Value = 0.0005207651961489113
$ws.Range("J20").Value = 0.0005230593158758034
$ws.Range("M20").Value = 1.5741895
$ws.Range("N20").Value = 3.148379
$ws.Range("O20").Value = 0.009700049718478087
$ws.Range("P20").Value = 0.006487676741301404
$ws.Range("Q20").Value = 0.4884022622418333
$ws.Range("R20").Value = 2.930413573451
$ws.Range("S20").Value = 0.000005051448294297433
$ws.Range("T20").Value = 0.000003393439757928474
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.3102563333333333
$ws.Range("H21").Value = 0.930769
$ws.Range("I21").Value = 0.0005207651961489113
$ws.Range("J21").Value = 0.0005230593158758034
$ws.Range("M21").Value = 72.07364666666666
$ws.Range("N21").Value = 216.22094
$ws.Range("O21").Value = 0.4441129584835175
$ws.Range("P21").Value = 0.4455535891391496
$ws.Range("Q21").Value = 22.36130534476222
$ws.Range("R21").Value = 201.25174810286
$ws.Range("S21").Value = 0.0002312785719369423
$ws.Range("T21").Value = 0.0002330509555211323
